$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.806.76'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -3.58%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.616.10'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -3.74%  '

$ws.Range('E4').Value = '  +0.26%  '

$ws.Range('B5').Value = 'USDC'
$ws.Range('C5').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.000'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.05%  '

$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '306.65'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.27%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3921'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.31%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3823'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.41%  '

$ws.Range('E9').Value = '  +0.17%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.365'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.15%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '49.20'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.61%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08424'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.88%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.86'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -6.11%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.026'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.48%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.544'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.48%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001275'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.62%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.621.14'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -3.01%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.18'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.99%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06906'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.84%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.03'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -6.17%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.814'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.03%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9997'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.04%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.39'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.39%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.802.21'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.57%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.449'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.51%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.863'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.60%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.12'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -4.79%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '156.88'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.21%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '139.20'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -5.16%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.243'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -10.68%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.862'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.39%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.485'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.37%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.784.74'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.72%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08023'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.67%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9766'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.48%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02876'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -7.11%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.571'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -6.00%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2663'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -5.53%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.09196'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.67%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.24'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.22%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.39'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.67%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.425'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.99%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7464'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -5.93%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.91'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.34%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6855'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.17%  '

$ws.Range('E46').Value = '  -4.62%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.057'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.80%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9999'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.04%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08257'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.66%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '133.13'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.39%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.204'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -9.90%  '
